$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.551.34"
$ws.Range("E2").Value = "  +2.39%  "

$ws.Range("D3").Value = "3.310.78"
$ws.Range("E3").Value = "  -0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.74"
$ws.Range("E5").Value = "  +4.01%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.92"
$ws.Range("E6").Value = "  -1.97%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.588"
$ws.Range("E8").Value = "  +2.52%  "

$ws.Range("D9").Value = "3.301.49"
$ws.Range("E9").Value = "  -0.34%  "

$ws.Range("E10").Value = "  +0.91%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.575"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "46.07"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("E13").Value = "  +3.02%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "634.60"
$ws.Range("E14").Value = "  +9.69%  "

$ws.Range("D15").Value = "3.843.04"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.42"
$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "67.691.51"
$ws.Range("E17").Value = "  +2.52%  "

$ws.Range("E18").Value = "  +1.55%  "

$ws.Range("D19").Value = "3.316.81"
$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.57"
$ws.Range("E20").Value = "  -0.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.85"
$ws.Range("E21").Value = "  -0.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.898"
$ws.Range("E22").Value = "  +0.36%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.51"
$ws.Range("E23").Value = "  -0.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.01"
$ws.Range("E24").Value = "  -0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "96.49"
$ws.Range("E25").Value = "  -1.13%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.99"
$ws.Range("E26").Value = "  +0.59%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.76"
$ws.Range("E27").Value = "  +2.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.53"
$ws.Range("E28").Value = "  +1.72%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "32.46"
$ws.Range("E29").Value = "  +6.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.53"
$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("E31").Value = "  -1.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "587.27"
$ws.Range("E32").Value = "  +4.63%  "

$ws.Range("D33").Value = "3.935.82"
$ws.Range("E33").Value = "  +4.30%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.89"
$ws.Range("E34").Value = "  +0.74%  "

$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  +0.80%  "

$ws.Range("B36").Value = "dogwifhat"
$ws.Range("C36").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.49"
$ws.Range("E36").Value = "  -5.41%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "55.48"
$ws.Range("E38").Value = "  -0.35%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.128"
$ws.Range("E39").Value = "  +1.18%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.23"
$ws.Range("E40").Value = "  +2.98%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.65"
$ws.Range("E41").Value = "  +2.73%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "32.38"
$ws.Range("E42").Value = "  -2.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.38"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").Value = "0.0₃0680"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  +0.73%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0412"
$ws.Range("E46").Value = "  +0.84%  "

$ws.Range("E47").Value = "  +1.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.01"
$ws.Range("E48").Value = "  +0.54%  "

$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.38"
$ws.Range("E49").Value = "  +12.16%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.53"
$ws.Range("E50").Value = "  +0.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "130.54"
$ws.Range("E51").Value = "  +2.89%  "

